$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# The old sheet was a Jan..Dec (2020/2021 header) month table in A:C.
# Remove the first five month rows (ENERO..MAYO) so the remaining
# JUNIO..DICIEMBRE block (old rows 7-13) slides up to rows 2-8 - this
# preserves the pre-existing per-cell text styling that lived on the
# old B7/B8 cells (now B2/B3).
# ---------------------------------------------------------------------
$ws.Range("A2:A6").EntireRow.Delete()

# Drop the now-unused column C (old "2021" column / stray monthly totals).
$ws.Columns.Item(3).Delete()

# ---- Header row: FECHA / TPD ----
$ws.Range("A1").Value = "FECHA"
$ws.Range("B1").Value = "TPD"

# ---- Column A becomes a date column (old month names replaced by dates) ----
$dates = @(43983,44013,44044,44075,44105,44136,44166,44197,44228,44256,44287,44317,44348)
for ($i = 0; $i -lt $dates.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $dates[$i]
}

# ---- New TPD values for the six months appended after DICIEMBRE ----
$ws.Range("B9").Value = 74
$ws.Range("B10").Value = 53
$ws.Range("B11").Value = 75
$ws.Range("B12").Value = 139
$ws.Range("B13").Value = 203
$ws.Range("B14").Value = 164

# ---- A handful of blank, pre-formatted rows left ready below the data ----
$ws.Range("A15").Value = ""
$ws.Range("A16").Value = ""
$ws.Range("A17").Value = ""
$ws.Range("A18").Value = ""
$ws.Range("A19").Value = ""

# Apply the date number format across the whole FECHA column in one pass,
# covering the header, every data row and the trailing blank rows.
$ws.Range("A1:A19").NumberFormat = "mm-dd-yy"

$ws.Range("D16").Select()
